# Added ifo GDP component analysis preprocessing:
# Fill in the next diagonal cell for each of the rows 11-20 in the
# matched-errors table (one new numeric value per row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K11").Value = 0.2970525035592049
$ws.Range("J12").Value = 0.2858677898194339
$ws.Range("I13").Value = 0.2775335613519331
$ws.Range("H14").Value = 0.2743085116504074
$ws.Range("G15").Value = 0.2534447081011285
$ws.Range("F16").Value = 0.2766837437271186
$ws.Range("E17").Value = 0.2867219094086165
$ws.Range("D18").Value = 0.1751453671933744
$ws.Range("C19").Value = 0.1965658720679752
$ws.Range("B20").Value = 0.4328090033804217
